$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G37").Value = "ultimates, gigas, base/meso zipper"
$ws.Range("G55").Value = "Havocs, Zippers, Shifters - up to macros"
$ws.Range("G57").Value = "all basics"
$ws.Range("G58").Value = "jammer, vessel up to macro, basic types"
$ws.Range("G59").Value = "macro tickers, engimas up to macro, booster"
$ws.Range("G60").Value = "carriers, outliers up to macro, zippers up to macro"
$ws.Range("G61").Value = "enigma up to macro, zipper up to macro, only base type"
$ws.Range("G63").Value = "disruptor, zippers, shifter (to macro)"
$ws.Range("G64").Value = "reflectors up to macro, gigas"
$ws.Range("G65").Value = "cores, armory, armored"
$ws.Range("G67").Value = "maintainer, protector, shields"
$ws.Range("G68").Value = "MacroVessel, MacroTeleporter"
$ws.Range("G75").Value = "All basic available, hyp havoc"
$ws.Range("G77").Value = "all shifters"
$ws.Range("G78").Value = "All basic available, hyp protector"
$ws.Range("G79").Value = "All basic available, hyp booster"
$ws.Range("G82").Value = "All basic available. Hyp maintainer"
$ws.Range("G83").Value = "All basic available, hyp armory"
$ws.Range("G85").Value = "All basic available, hyp jammer"
$ws.Range("G86").Value = "All basic available, hyp disruptor"
$ws.Range("G87").Value = "All basic available, hyp disruptor & jammer"
$ws.Range("G88").Value = "carrier, core, basic types up to kilo"
$ws.Range("G90").Value = "tickers, enigmas"
$ws.Range("G93").Value = "All basic available, hyp protector, armory, maintainer, disruptor"
$ws.Range("G95").Value = "All basic available, hyp havoc, booster, jammer"
$ws.Range("G97").Value = "All basic available, all buff/debuff"
$ws.Range("G98").Value = "All basic available, all the non hyper buffers/debuffer"

Write-Output "done"
